# Making changes in Search module
# Adds TestCase_B25..B30 ("type ahead" / autocomplete related cases) to the
# "Test Cases" sheet, fixes the Results value of TestCase_B6 (row 11, was
# incorrectly PASS) to SKIP, and extends the sheet's used range/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 11 (TestCase_B6): Results column was wrongly "PASS" -> fix to "SKIP"
$ws.Range("E11").Value = "SKIP"

# --- Row 13 (TestCase_B8): Description re-pointed to the "ALL content type
#     count" text (already present elsewhere in the workbook).
$ws.Range("C13").Value = "Verify that ALL content type count is equal to the sum of the counts of other content types"

# Helper: clone formatting (borders/fill/alignment) from a template cell that
# already carries the visual style we need, then write the real value. This
# reuses existing style entries instead of fabricating new ones.
function Set-FormattedCell($ws, $addr, $templateAddr, $value) {
    $ws.Range($templateAddr).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $value
}

# Template cells already on the sheet:
#   D2  -> plain bordered cell (A/D columns style)
#   B2  -> bordered + shaded + wrap-text cell (B column style)
#   E2  -> plain bordered, non-wrapped cell (C/E columns, no-wrap style)
$plain = "D2"
$shadedWrap = "B2"
$plainNoWrap = "E2"

# --- Row 26: TestCase_B25
Set-FormattedCell $ws "A26" $plain "TestCase_B25"
Set-FormattedCell $ws "B26" $shadedWrap "TBD-1"
Set-FormattedCell $ws "C26" $plainNoWrap "Verify that autocomplete feature is working correctly"
Set-FormattedCell $ws "D26" $plain "Y"
Set-FormattedCell $ws "E26" $plainNoWrap "SKIP"

# --- Row 27: TestCase_B26 (tall wrapped row)
Set-FormattedCell $ws "A27" $plain "TestCase_B26"
Set-FormattedCell $ws "B27" $shadedWrap "TBD-2"
Set-FormattedCell $ws "C27" $plainNoWrap "Verify that following sections get displayed in type ahead:`na)Categories`nb)Articles`nc)Patents`nd)People"
$ws.Range("C27").WrapText = $true
Set-FormattedCell $ws "D27" $plain "Y"
Set-FormattedCell $ws "E27" $plainNoWrap "SKIP"
$ws.Rows.Item(27).RowHeight = 75

# --- Row 28: TestCase_B27
Set-FormattedCell $ws "A28" $plain "TestCase_B27"
Set-FormattedCell $ws "B28" $shadedWrap "TBD-3"
Set-FormattedCell $ws "C28" $plainNoWrap "Verify that 4 suggested categories get displayed in type ahead and the typed keyword is present in all the 4 categories"
Set-FormattedCell $ws "D28" $plain "Y"
Set-FormattedCell $ws "E28" $plainNoWrap "SKIP"

# --- Row 29: TestCase_B28
Set-FormattedCell $ws "A29" $plain "TestCase_B28"
Set-FormattedCell $ws "B29" $shadedWrap "TBD-4"
Set-FormattedCell $ws "C29" $plainNoWrap "Verify that 4 suggested articles get displayed in type ahead and the typed keyword is present in all the 4 articles"
Set-FormattedCell $ws "D29" $plain "Y"
Set-FormattedCell $ws "E29" $plainNoWrap "SKIP"

# --- Row 30: TestCase_B29
Set-FormattedCell $ws "A30" $plain "TestCase_B29"
Set-FormattedCell $ws "B30" $shadedWrap "TBD-5"
Set-FormattedCell $ws "C30" $plainNoWrap "Verify that 4 suggested patents get displayed in type ahead and the typed keyword is present in all the 4 patents"
Set-FormattedCell $ws "D30" $plain "Y"
Set-FormattedCell $ws "E30" $plainNoWrap "SKIP"

# --- Row 31: TestCase_B30
Set-FormattedCell $ws "A31" $plain "TestCase_B30"
Set-FormattedCell $ws "B31" $shadedWrap "TBD-6"
Set-FormattedCell $ws "C31" $plainNoWrap "Verify that 4 suggested people get displayed in type ahead and the typed keyword is present in all the 4 people"
Set-FormattedCell $ws "D31" $plain "Y"
Set-FormattedCell $ws "E31" $plainNoWrap "PASS"

# --- Extend the visible selection to cover the new rows (dimension is
#     recomputed automatically from the used range).
$ws.Range("D2:D31").Select() | Out-Null
